# Update example: add 2 outputs of one process
# The source data had an output tag "O6" that should become "O7".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3/B3 currently hold "O6" -> rename to "O7"
$ws.Range("A3").Value = "O7"
$ws.Range("B3").Value = "O7"

# Update the active selection to F3, matching the saved view state
$ws.Range("F3").Select() | Out-Null
